$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Column B/C (coin name / link) swaps ---
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"

# --- Column D (price) updates ---
Set-TextValue $ws.Range("D2") "27.400.56"
Set-TextValue $ws.Range("D3") "1.863.26"
Set-TextValue $ws.Range("D5") "315.20"
Set-TextValue $ws.Range("D6") "1.004"
Set-TextValue $ws.Range("D7") "0.4647"
Set-TextValue $ws.Range("D8") "0.3717"
Set-TextValue $ws.Range("D9") "0.07358"
Set-TextValue $ws.Range("D10") "0.8864"
Set-TextValue $ws.Range("D11") "0.07919"
Set-TextValue $ws.Range("D13") "1.849.79"
Set-TextValue $ws.Range("D14") "5.408"
Set-TextValue $ws.Range("D15") "6.593"
Set-TextValue $ws.Range("D16") "92.33"
Set-TextValue $ws.Range("D17") "1.005"
Set-TextValue $ws.Range("D18") "0.000008887"
Set-TextValue $ws.Range("D20") "14.86"
Set-TextValue $ws.Range("D21") "27.448.29"
Set-TextValue $ws.Range("D23") "10.55"
Set-TextValue $ws.Range("D24") "2.137.56"
Set-TextValue $ws.Range("D25") "1.896"
Set-TextValue $ws.Range("D26") "153.02"
Set-TextValue $ws.Range("D27") "18.44"
Set-TextValue $ws.Range("D28") "2.074"
Set-TextValue $ws.Range("D29") "5.137"
Set-TextValue $ws.Range("D30") "116.54"
Set-TextValue $ws.Range("D31") "0.08897"
Set-TextValue $ws.Range("D32") "0.7567"
Set-TextValue $ws.Range("D33") "3.025"
Set-TextValue $ws.Range("D34") "1.166"
Set-TextValue $ws.Range("D35") "4.498"
Set-TextValue $ws.Range("D37") "0.01964"
Set-TextValue $ws.Range("D38") "1.079"
Set-TextValue $ws.Range("D39") "2.988"
Set-TextValue $ws.Range("D40") "0.05255"
Set-TextValue $ws.Range("D41") "7.148"
Set-TextValue $ws.Range("D42") "0.5176"
Set-TextValue $ws.Range("D43") "0.1645"
Set-TextValue $ws.Range("D44") "8.358"
Set-TextValue $ws.Range("D45") "0.4853"
Set-TextValue $ws.Range("D46") "10.36"
Set-TextValue $ws.Range("D47") "1.004"
Set-TextValue $ws.Range("D48") "104.05"
Set-TextValue $ws.Range("D50") "0.06249"
Set-TextValue $ws.Range("D51") "65.92"

# --- Column E (volume/percent change) updates ---
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  +5.11%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("E36").Value = "  +10.10%  "
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  +2.34%  "
